# Updated MCH102 to MCH251: add a new descriptive row (row 2) with the
# identifier/levelOfDescription/extentAndMedium/notes for MCH237-1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Values for the new row --------------------------------------------
$ws.Range("A2").Value = "MCH237-1"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 2C | GRAP COUNT NUMER: NONE"

# -- Formatting: Calibri 10pt, automatic (theme) text colour -----------
# Applied per contiguous block since non-contiguous (union) ranges only
# style their first area.
foreach ($addr in @("A2", "C2:E2", "G2:H2")) {
    $rng = $ws.Range($addr)
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 10
    $rng.Font.ThemeColor = 1
}

$f2 = $ws.Range("F2")
$f2.Font.Name = "Calibri"
$f2.Font.Size = 10
$f2.Font.ThemeColor = 1

# -- Selection / frozen header row restated on the new active row ------
$ws.Range("A2:H2").Select()
$excel.ActiveWindow.FreezePanes = $true
